# Generate Report for Handback
#
# - The "Status" value that used to read "Ready for handoff" (for the
#   c247b7ae-4619-441b-87f5-07da20779c83 file) now reads
#   "Handback transform failed". That text is shared by the Overview
#   sheet (columns E/F) and by the per-locale "Status" column (C) on
#   both the zh-cn and de-de sheets, so all four cells are updated.
# - The zh-cn and de-de sheets each record an "Error Detail" message
#   in column P for that same row, explaining the handback/handoff
#   filename mismatch, and column P is widened so the message is
#   readable.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

$newStatus = "Handback transform failed"
$overview.Range("E3").Value = $newStatus
$overview.Range("F3").Value = $newStatus
$zhcn.Range("C3").Value = $newStatus
$dede.Range("C3").Value = $newStatus

# Record the handback/handoff file name mismatch error for each locale.
$zhcn.Range("P3").Value = "Handback file name: v403zv40.sms is different with handoff file name: c247b7ae-4619-441b-87f5-07da20779c83.b81ba540e37c5d7d19079baf568a260122e9752a.zh-cn."
$dede.Range("P3").Value = "Handback file name: v403zv40.sms is different with handoff file name: c247b7ae-4619-441b-87f5-07da20779c83.b81ba540e37c5d7d19079baf568a260122e9752a.de-de."

# Widen the "Error Detail" column (P) on both sheets so the message fits.
# ColumnWidth is expressed in characters; the saved OOXML column width
# adds a fixed 5/6 character padding on top of it, so back that out here
# to land on an OOXML width of exactly 40.
$zhcn.Columns.Item(16).ColumnWidth = 40 - 5/6
$dede.Columns.Item(16).ColumnWidth = 40 - 5/6
